$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.020705105263572
$ws.Cells.Item(2, 4).Value = 1.026649187597143
$ws.Cells.Item(2, 5).Value = 1.021655281550852
$ws.Cells.Item(2, 6).Value = 1.033991353004535
$ws.Cells.Item(2, 9).Value = 1.030318399236568
$ws.Cells.Item(2, 10).Value = 1.025900707851144
$ws.Cells.Item(2, 11).Value = 1.029471217959963
$ws.Cells.Item(2, 12).Value = 1.024491970833195
$ws.Cells.Item(2, 13).Value = 1.036792110909094
$ws.Cells.Item(2, 14).Value = 1.012624279800934
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.021903196879089
$ws.Cells.Item(3, 4).Value = 1.027576811208665
$ws.Cells.Item(3, 5).Value = 1.022679554830699
$ws.Cells.Item(3, 6).Value = 1.036285496694253
$ws.Cells.Item(3, 9).Value = 1.030670969569198
$ws.Cells.Item(3, 10).Value = 1.026734806261003
$ws.Cells.Item(3, 11).Value = 1.030206075272806
$ws.Cells.Item(3, 12).Value = 1.025322147436535
$ws.Cells.Item(3, 13).Value = 1.038891395419815
$ws.Cells.Item(3, 14).Value = 1.012907789684852
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.02267761201774
$ws.Cells.Item(4, 4).Value = 1.028175921504303
$ws.Cells.Item(4, 5).Value = 1.023341936307388
$ws.Cells.Item(4, 6).Value = 1.037763094351598
$ws.Cells.Item(4, 9).Value = 1.030896520502329
$ws.Cells.Item(4, 10).Value = 1.027273228892184
$ws.Cells.Item(4, 11).Value = 1.030679834863818
$ws.Cells.Item(4, 12).Value = 1.025858339136397
$ws.Cells.Item(4, 13).Value = 1.040242608526113
$ws.Cells.Item(4, 14).Value = 1.013090606757572
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.023002981386308
$ws.Cells.Item(5, 4).Value = 1.028427521273367
$ws.Cells.Item(5, 5).Value = 1.023620310347898
$ws.Cells.Item(5, 6).Value = 1.038382670672587
$ws.Cells.Item(5, 9).Value = 1.030990726257877
$ws.Cells.Item(5, 10).Value = 1.027499274665899
$ws.Cells.Item(5, 11).Value = 1.030878589210967
$ws.Cells.Item(5, 12).Value = 1.026083520162919
$ws.Cells.Item(5, 13).Value = 1.04080897900019
$ws.Cells.Item(5, 14).Value = 1.013167312481795
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.023057601030056
$ws.Cells.Item(6, 4).Value = 1.028469750456049
$ws.Cells.Item(6, 5).Value = 1.02366704530708
$ws.Cells.Item(6, 6).Value = 1.038486607200316
$ws.Cells.Item(6, 9).Value = 1.031006507795162
$ws.Cells.Item(6, 10).Value = 1.02753721082036
$ws.Cells.Item(6, 11).Value = 1.03091193673598
$ws.Cells.Item(6, 12).Value = 1.026121315402708
$ws.Cells.Item(6, 13).Value = 1.040903977632675
$ws.Cells.Item(6, 14).Value = 1.013180182907087
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.022681960378629
$ws.Cells.Item(7, 4).Value = 1.028179284434431
$ws.Cells.Item(7, 5).Value = 1.023345656309085
$ws.Cells.Item(7, 6).Value = 1.037771379418999
$ws.Cells.Item(7, 9).Value = 1.030897781699754
$ws.Cells.Item(7, 10).Value = 1.02727625052939
$ws.Cells.Item(7, 11).Value = 1.030682492252404
$ws.Cells.Item(7, 12).Value = 1.025861348931838
$ws.Cells.Item(7, 13).Value = 1.040250182947668
$ws.Cells.Item(7, 14).Value = 1.013091632293516
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.021110180269121
$ws.Cells.Item(8, 4).Value = 1.026962916899311
$ws.Cells.Item(8, 5).Value = 1.022001522251602
$ws.Cells.Item(8, 6).Value = 1.034768114547154
$ws.Cells.Item(8, 9).Value = 1.030438089320833
$ws.Cells.Item(8, 10).Value = 1.026182865160489
$ws.Cells.Item(8, 11).Value = 1.02971992893477
$ws.Cells.Item(8, 12).Value = 1.024772739056566
$ws.Cells.Item(8, 13).Value = 1.037503078975409
$ws.Cells.Item(8, 14).Value = 1.012720224870104
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.01833396256332
$ws.Cells.Item(9, 4).Value = 1.024810777511521
$ws.Cells.Item(9, 5).Value = 1.019629863701291
$ws.Cells.Item(9, 6).Value = 1.02942165516363
$ws.Cells.Item(9, 9).Value = 1.02960810586051
$ws.Cells.Item(9, 10).Value = 1.024246132596182
$ws.Cells.Item(9, 11).Value = 1.028010293275922
$ws.Cells.Item(9, 12).Value = 1.022846784315859
$ws.Cells.Item(9, 13).Value = 1.032605842744724
$ws.Cells.Item(9, 14).Value = 1.012060869478796
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.016478495162028
$ws.Cells.Item(10, 4).Value = 1.023369945803
$ws.Cells.Item(10, 5).Value = 1.018046484087756
$ws.Cells.Item(10, 6).Value = 1.025818378529811
$ws.Cells.Item(10, 9).Value = 1.02904116740109
$ws.Cells.Item(10, 10).Value = 1.022948034490629
$ws.Cells.Item(10, 11).Value = 1.026861285435581
$ws.Cells.Item(10, 12).Value = 1.021557479212489
$ws.Cells.Item(10, 13).Value = 1.029300753666374
$ws.Cells.Item(10, 14).Value = 1.011617952937474
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.015673894273775
$ws.Cells.Item(11, 4).Value = 1.022744565069668
$ws.Cells.Item(11, 5).Value = 1.017360283769074
$ws.Cells.Item(11, 6).Value = 1.024248325050152
$ws.Cells.Item(11, 9).Value = 1.028792400365991
$ws.Cells.Item(11, 10).Value = 1.022384253912313
$ws.Cells.Item(11, 11).Value = 1.026361514502574
$ws.Cells.Item(11, 12).Value = 1.020997891789603
$ws.Cells.Item(11, 13).Value = 1.027859547825072
$ws.Cells.Item(11, 14).Value = 1.011425357721418
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.015374848517087
$ws.Cells.Item(12, 4).Value = 1.02251204282278
$ws.Cells.Item(12, 5).Value = 1.017105306722564
$ws.Cells.Item(12, 6).Value = 1.023663615330607
$ws.Cells.Item(12, 9).Value = 1.028699500432461
$ws.Cells.Item(12, 10).Value = 1.022174582142023
$ws.Cells.Item(12, 11).Value = 1.02617553644528
$ws.Cells.Item(12, 12).Value = 1.020789835877527
$ws.Cells.Item(12, 13).Value = 1.027322660486506
$ws.Cells.Item(12, 14).Value = 1.011353696440033
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.015439003129124
$ws.Cells.Item(13, 4).Value = 1.022561930016088
$ws.Cells.Item(13, 5).Value = 1.017160004372097
$ws.Cells.Item(13, 6).Value = 1.023789107201511
$ws.Cells.Item(13, 9).Value = 1.028719450331522
$ws.Cells.Item(13, 10).Value = 1.022219569211046
$ws.Cells.Item(13, 11).Value = 1.026215444858862
$ws.Cells.Item(13, 12).Value = 1.020834473693552
$ws.Cells.Item(13, 13).Value = 1.02743789593457
$ws.Cells.Item(13, 14).Value = 1.011369073607098
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.015649178784184
$ws.Cells.Item(14, 4).Value = 1.022725349382569
$ws.Cells.Item(14, 5).Value = 1.017339209178348
$ws.Cells.Item(14, 6).Value = 1.024200024090492
$ws.Cells.Item(14, 9).Value = 1.02878473138702
$ws.Cells.Item(14, 10).Value = 1.022366927675308
$ws.Cells.Item(14, 11).Value = 1.026346148475036
$ws.Cells.Item(14, 12).Value = 1.020980697927719
$ws.Cells.Item(14, 13).Value = 1.027815200622363
$ws.Cells.Item(14, 14).Value = 1.01141943669072
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.015778650779648
$ws.Cells.Item(15, 4).Value = 1.022826007115593
$ws.Cells.Item(15, 5).Value = 1.017449610915355
$ws.Cells.Item(15, 6).Value = 1.024453000423633
$ws.Cells.Item(15, 9).Value = 1.028824887245862
$ws.Cells.Item(15, 10).Value = 1.022457685711148
$ws.Cells.Item(15, 11).Value = 1.026426633998766
$ws.Cells.Item(15, 12).Value = 1.021070764868284
$ws.Cells.Item(15, 13).Value = 1.028047462511792
$ws.Cells.Item(15, 14).Value = 1.011450450733558
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.016531868754318
$ws.Cells.Item(16, 4).Value = 1.02341141850924
$ws.Cells.Item(16, 5).Value = 1.018092012281884
$ws.Cells.Item(16, 6).Value = 1.025922367029785
$ws.Cells.Item(16, 9).Value = 1.029057607813501
$ws.Cells.Item(16, 10).Value = 1.022985414702752
$ws.Cells.Item(16, 11).Value = 1.026894405998582
$ws.Cells.Item(16, 12).Value = 1.02159458929802
$ws.Cells.Item(16, 13).Value = 1.029396185559603
$ws.Cells.Item(16, 14).Value = 1.01163071768895
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.017004025265852
$ws.Cells.Item(17, 4).Value = 1.023778229617478
$ws.Cells.Item(17, 5).Value = 1.018494814155665
$ws.Cells.Item(17, 6).Value = 1.026841401640496
$ws.Cells.Item(17, 9).Value = 1.029202706503062
$ws.Cells.Item(17, 10).Value = 1.023315988532681
$ws.Cells.Item(17, 11).Value = 1.027187223780863
$ws.Cells.Item(17, 12).Value = 1.021922817437518
$ws.Cells.Item(17, 13).Value = 1.030239473679205
$ws.Cells.Item(17, 14).Value = 1.011743576814379
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.017279313516683
$ws.Cells.Item(18, 4).Value = 1.023992040735529
$ws.Cells.Item(18, 5).Value = 1.018729705331257
$ws.Cells.Item(18, 6).Value = 1.027376515654036
$ws.Cells.Item(18, 9).Value = 1.029287024029014
$ws.Cells.Item(18, 10).Value = 1.023508643323138
$ws.Cells.Item(18, 11).Value = 1.027357803294642
$ws.Cells.Item(18, 12).Value = 1.022114141217362
$ws.Cells.Item(18, 13).Value = 1.030730379825314
$ws.Cells.Item(18, 14).Value = 1.011809327619309
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.01737316074371
$ws.Cells.Item(19, 4).Value = 1.024064920563918
$ws.Cells.Item(19, 5).Value = 1.018809787702368
$ws.Cells.Item(19, 6).Value = 1.027558817137908
$ws.Cells.Item(19, 9).Value = 1.029315720637294
$ws.Cells.Item(19, 10).Value = 1.023574306038636
$ws.Cells.Item(19, 11).Value = 1.027415929955435
$ws.Cells.Item(19, 12).Value = 1.022179356389397
$ws.Cells.Item(19, 13).Value = 1.030897603027523
$ws.Cells.Item(19, 14).Value = 1.011831733749115
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.016953379030581
$ws.Cells.Item(20, 4).Value = 1.023738889138539
$ws.Cells.Item(20, 5).Value = 1.018451603180413
$ws.Cells.Item(20, 6).Value = 1.026742895841521
$ws.Cells.Item(20, 9).Value = 1.02918717151093
$ws.Cells.Item(20, 10).Value = 1.023280538011105
$ws.Cells.Item(20, 11).Value = 1.027155829589059
$ws.Cells.Item(20, 12).Value = 1.021887614736743
$ws.Cells.Item(20, 13).Value = 1.030149097451547
$ws.Cells.Item(20, 14).Value = 1.011731476176558
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.015587292335274
$ws.Cells.Item(21, 4).Value = 1.022677232778429
$ws.Cells.Item(21, 5).Value = 1.017286440389754
$ws.Cells.Item(21, 6).Value = 1.024079061711569
$ws.Cells.Item(21, 9).Value = 1.028765521490623
$ws.Cells.Item(21, 10).Value = 1.022323541448996
$ws.Cells.Item(21, 11).Value = 1.026307668959441
$ws.Cells.Item(21, 12).Value = 1.020937644090503
$ws.Cells.Item(21, 13).Value = 1.027704137205098
$ws.Cells.Item(21, 14).Value = 1.011404609418946
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.014727328638965
$ws.Cells.Item(22, 4).Value = 1.02200840613832
$ws.Cells.Item(22, 5).Value = 1.016553323961247
$ws.Cells.Item(22, 6).Value = 1.022395373499426
$ws.Cells.Item(22, 9).Value = 1.028497537050144
$ws.Cells.Item(22, 10).Value = 1.021720341830101
$ws.Cells.Item(22, 11).Value = 1.025772422492601
$ws.Cells.Item(22, 12).Value = 1.020339199638171
$ws.Cells.Item(22, 13).Value = 1.026157849076502
$ws.Cells.Item(22, 14).Value = 1.011198383957277
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.015183313018073
$ws.Cells.Item(23, 4).Value = 1.022363090372079
$ws.Cells.Item(23, 5).Value = 1.016942014520878
$ws.Cells.Item(23, 6).Value = 1.023288781711939
$ws.Cells.Item(23, 9).Value = 1.028639874725279
$ws.Cells.Item(23, 10).Value = 1.022040252638051
$ws.Cells.Item(23, 11).Value = 1.026056355226079
$ws.Cells.Item(23, 12).Value = 1.02065655750343
$ws.Cells.Item(23, 13).Value = 1.02697843818906
$ws.Cells.Item(23, 14).Value = 1.01130777582254
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.016976264232682
$ws.Cells.Item(24, 4).Value = 1.023756665851562
$ws.Cells.Item(24, 5).Value = 1.018471128534563
$ws.Cells.Item(24, 6).Value = 1.026787409284736
$ws.Cells.Item(24, 9).Value = 1.029194192081861
$ws.Cells.Item(24, 10).Value = 1.023296557080785
$ws.Cells.Item(24, 11).Value = 1.027170015941063
$ws.Cells.Item(24, 12).Value = 1.021903521711966
$ws.Cells.Item(24, 13).Value = 1.030189937572576
$ws.Cells.Item(24, 14).Value = 1.011736944175078
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.019052482822012
$ws.Cells.Item(25, 4).Value = 1.025368213001057
$ws.Cells.Item(25, 5).Value = 1.020243382105151
$ws.Cells.Item(25, 6).Value = 1.030810529807987
$ws.Cells.Item(25, 9).Value = 1.029825061129775
$ws.Cells.Item(25, 10).Value = 1.024748034211153
$ws.Cells.Item(25, 11).Value = 1.028453891410406
$ws.Cells.Item(25, 12).Value = 1.023345617648376
$ws.Cells.Item(25, 13).Value = 1.03387882258319
$ws.Cells.Item(25, 14).Value = 1.01223191400704
